$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..55 (generation, elite_count)
$data = @(
    @(0, 2),
    @(1, 1),
    @(5, 1),
    @(6, 4),
    @(7, 1),
    @(8, 2),
    @(14, 1),
    @(16, 1),
    @(18, 1),
    @(20, 1),
    @(21, 1),
    @(22, 3),
    @(26, 1),
    @(27, 1),
    @(31, 1),
    @(32, 1),
    @(33, 2),
    @(35, 2),
    @(36, 2),
    @(37, 2),
    @(39, 1),
    @(41, 1),
    @(42, 2),
    @(43, 1),
    @(45, 2),
    @(48, 1),
    @(50, 2),
    @(51, 6),
    @(52, 1),
    @(53, 2),
    @(54, 3),
    @(55, 1),
    @(56, 1),
    @(57, 2),
    @(60, 1),
    @(61, 1),
    @(62, 2),
    @(63, 2),
    @(65, 1),
    @(66, 3),
    @(67, 8),
    @(68, 2),
    @(69, 1),
    @(71, 1),
    @(73, 2),
    @(74, 1),
    @(77, 1),
    @(78, 2),
    @(79, 5),
    @(81, 1),
    @(82, 1),
    @(83, 2),
    @(84, 3),
    @(86, 3),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-obsolete rows 56..61 that existed in the original sheet
$lastRow = $data.Count + 2
$ws.Range("A" + $lastRow + ":B61").EntireRow.Delete()

Write-Host "Updated elite_count_by_generation data to " $data.Count " rows"
